$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 39
$srcRow = 38

# Copy the date cell's style/number-format from the row above, then set the new date value
$ws.Cells.Item($newRow, 1).NumberFormat = $ws.Cells.Item($srcRow, 1).NumberFormat
$ws.Cells.Item($newRow, 1).Value = 45887

$ws.Cells.Item($newRow, 2).Value = 6449.14990234375
$ws.Cells.Item($newRow, 3).Value = 6455.35009765625
$ws.Cells.Item($newRow, 4).Value = 6437.7001953125
$ws.Cells.Item($newRow, 5).Value = 6445.02001953125
$ws.Cells.Item($newRow, 6).Value = 4162290000
$ws.Cells.Item($newRow, 7).Value = -0.0001007631807854
$ws.Cells.Item($newRow, 8).Value = -1

$ws.Cells.Item($newRow, 9).Value = 'strategy mstr keep underperforming bitcoin whats driving the divergence online casino company acquires spotlightvegas medium advisory federal government to make housing announcement in stratford cold wallet among the highest trending crypto with m raised while xrp stay bullish dot price eye colorado springsbased auria win million contract watch live trump putin share warm greeting handshake a alaska summit kick off rddt deadline monday rosen a top ranked law firm encourages reddit inc investor with air canada could shut down completely unless the airline and it flight attendant reach a deal air canada could shut down completely unless the airline and it flight attendant reach a deal advocate mark passage of ct towing reform law jetblue add more florida flight msnbc mouthpiece claim climate change is more dangerous than crime in dc caitong international asset management co ltd acquires share of visa inc nysev rddt deadline monday rosen a top ranked law firm encourages reddit inc investor with loss in excess of k to secure counsel before important august deadline in security class action rddt rosen recognized investor counsel encourages biohaven ltd investor to secure counsel before target and ulta mutually agree not to renew partnership launched in rosen recognized investor counsel encourages biohaven ltd investor to secure counsel before important deadline in security class action bhvn state where resident could see the biggest tax benefit from trump big beautiful bill trump is open to doing business with russia that could be tricky tariff sale slump hit deere but analyst see greener pasture ahead duke energy seek to extend operation for another year at bad creek supporting unprecedented growth in the carolina inventhelp inventor develops new laserbased tool for plumber tro government of canada provides early decision on york factory first nation ten shilling aerodrome project in manitoba saks say hyperpersonalized homepage boost revenue per visitor by monthly u home payment fall to a month low report trump say no new usrussia trade deal until war is settledahead of putin meeting live update job coming to dearborn county federal grand jury indicts new orleans mayor after long corruption probe little pepe lilpepe price forecast why is not an unrealistic longterm target gamblingcom analyst slash their forecast following q result forward industry recovers after q revenue decline and earnings miss halle berry is busy cooking cleaning and mothering in a picture perfect response to ex what is an upsidedown car loan and how do you get out of it breaking down sandisk analyst share their view what analyst rating have to say about restaurant brand intl assessing netstreit insight from financial analyst unpacking the latest option trading trend in centrus energy summer travel isnt a easy a it used to be for airline these analyst revise their forecast on deere after q result spotlight on conocophillips analyzing the surge in option activity demystifying amcor insight from analyst review what analyst are saying about nnn reit stock farleys model t moment bold ev strategy in a risky policy era first merchant corporation announces cash dividend pe ratio insight for targa resource what analyst are saying about broadstone net lease stock inside the backyard housing boom in california north river marketplace celebrates five year of selling food to resident jollibee and global girl group katseye ignite partnership to create joyful new connection with their north american fan m presale and blockdag academy put stellar and bittensor in the rearview for the best crypto to buy now sunrun run stock soar a well fargo nearly double price target liverpool sign teenage italian defender giovanni leoni liverpool sign teenage italian defender giovanni leoni cfo must take ownership of managing culture risk auditboard these analyst increase their forecast on advance auto part following upbeat q earnings novo nordisks lead fade eli lilly hold course for outsized gain say analyst after nearly year a legendary japanese sport car return with an unexpected twist gear patrol fact focus no taxpayer will not receive new stimulus check this summer fact focus no taxpayer will not receive new stimulus check this summer fact focus no taxpayer will not receive new stimulus check this summer fact focus no taxpayer will not receive new stimulus check this summer fact focus no taxpayer will not receive new stimulus check this summer got microplastic pollution anxiety chris hemsworth ha a cure these analyst revise their forecast on applied industrial technology following q earnings layer brett against pepe forecast the fight for ethereums leading meme spot begin just of californian can afford a home v nationally a strike loom air canada flight attendant reject call for arbitration the new york time air canada flight attendant reject call for arbitration unitedhealths best day since intels best weekly gain in decade whats moving market friday special livery package for gt celebrates th anniversary of andial is trump america ceo dogecoin news today m dogecoin whale transfer to binance trigger volatility concern senator call for probe into meta after news report on ai conversation with child sharplink share drop after q revenue decline and large noncash eth impairment retail sale rose in july a some shopper stepped up purchase ahead of tariff trump want to make intel great again what that mean for tsmc fed end specialized crypto bank oversight folding it into normal supervision these analyst boost their forecast on dillards following strong q earnings seven new toll brother model home open at the station in sunnyvale california burger bending your budget blame it on climate change best memecoin to buy now pepeto the next bnb of memecoins could turn investment into million agco order immediate suspension of north bay max multi mart lottery and alcohol licence in response to human trafficking charge washington dc trump administration court hearing for police takeover lawsuit stanley firstever warehouse sale ha deal up to off employee at this longtime newport news travel agency have over year of experience combined industrials stock whale activity in today session popular bakery hoping to avoid closure after zoning complaint samsung is giving away a tb ssd with this odyssey g monitor deepwater asset management llc ha million stock holding in broadcom inc nasdaqavgo billionaire family office join berkshire unitedhealth bet yahoo finance marlton partner file preliminary proxy statement related to election of director for the degree capital board of director lego newest tuxedo cat is here and it the purrrfect addition to any cat lover home jeff bezos lauren sanchez suffer a rocky start to their marriage with a devastating update yahoocom jeff bezos lauren sanchez suffer a rocky start to their marriage with a devastating update sheknows how ai and human expertise are changing wealth management service target and ulta beauty are ending their store partnership target and ulta beauty are ending their store partnership solana whale activity resume stellars rally why cold wallet is going viral with m raised in week sol price swing push solana investor toward layer brett and it staking reward top coding apps in and how to pick the best one perrone unlike at columbia trump attack on ucla is aimed at taxpayer money uhoh trump ha lost his rizz among the youth that voted for him republican make pitch for statewide office senate illinois dems say state should be blueprint in national fight against trump texas gop chairman call illinois a winnable purple state u house minority leader wa keynote speaker at annual democratic county chair brunch editorial hiding the homeless wont solve the problem shribman back to the future with pat buchanan letter trump administration mistreats federal worker federal grand jury indicts new orleans mayor latoya cantrell after long corruption probe new orleans mayor indicted over allegation of trying to hide relationship with bodyguard federal grand jury indicts new orleans mayor latoya cantrell after long corruption probe federal grand jury indicts new orleans mayor latoya cantrell after long corruption probe federal grand jury indicts new orleans mayor latoya cantrell after long corruption probe federal grand jury indicts new orleans mayor latoya cantrell after long corruption probe federal grand jury indicts new orleans mayor latoya cantrell after long corruption probe new orleans mayor latoya cantrell indicted by federal grand jury following lengthy corruption investigation texas gov abbott call for immediate second special session on redistricting bad diet too many med no exercise a look inside the latest make america healthy again report trump arrives at alaska military base for summit with putin on the russiaukraine war trump arrives at alaska military base for summit with putin on the russiaukraine war trump arrives at alaska military base for summit with putin on the russiaukraine war bad diet too many med no exercise a look inside the latest make america healthy again report trump is open to doing business with russia that could be tricky antivaxxers are making excuse for the cdc shooter putin tease new nuke deal ahead of alaska meeting with trump talking past each other why the useu dispute over free speech is set to escalate trump say no new usrussia trade deal until war is settledahead of putin meeting live update live now trump arrives in anchorage alaska federal grand jury indicts new orleans mayor after long corruption probe what adyens tariff trouble portend for payment company marjorie taylor greene keep buying unitedhealth stock along with fellow investor buffett and burry new orleans mayor indicted over allegation of trying to hide relationship with bodyguard halle berry is busy cooking cleaning and mothering in a picture perfect response to ex audience erupts at rep jasmine crockett town hall amid verbal attack audience erupts at rep jasmine crockett town hall amid verbal attack in america national park are more than scenic theyre sacred but they were created at a cost to native american cornyn paxton in dead heat in u senate primary latest poll show ice plan to open a many a three new detention center in rural colorado report say ice plan to open a many a three new detention center in rural colorado report say automatic gratuity may not be eligible for new trump tax break on tip how ai is changing the music business switzerland scale back f order following tariff threat republican view on supporting ukraine swing dramatically ahead of trump meeting with putin new poll say kiss criticized trump prior to kennedy center honor fact focus no taxpayer will not receive new stimulus check this summer fact focus no taxpayer will not receive new stimulus check this summer fact focus no taxpayer will not receive new stimulus check this summer fact focus no taxpayer will not receive new stimulus check this summer fact focus no taxpayer will not receive new stimulus check this summer fact focus no taxpayer will not receive new stimulus check this summer fact focus no taxpayer will not receive new stimulus check this summer fact focus no taxpayer will not receive new stimulus check this summer white house fire back a trump use of dc police for crime crackdown is placed in the crosshairs by dems witnessing deportation is devastating for child dhs say noem is temporarily staying in military housing after facing vicious doxxing death threat texas house speaker expects enough democrat to be back for quorum texas house speaker expects enough democrat to be back for quorum texas house speaker expects enough democrat to be back for quorum texas house speaker expects enough democrat to be back for quorum texas house speaker expects enough democrat to be back for quorum texas house democrat accept skokie official invite to north suburban mosque ap news summary at pm edt la nina weather return to california what to know trump just launched a new assault on dc home rule it brazenly illegal consumer group urge fhfa not to mix crypto with mortgage washington dc trump administration court hearing for police takeover lawsuit newsom say he will end the trump presidency reveals plan to bypass ca redistricting commission trump and putin pivotal alaska summit could reshape the war in ukraine trump mull lower refugee cap of with most space for white south african hawley open probe into meta after report of ai romantic exchange with minor septum urge rider to start planning now for upcoming service cut that begin on aug weekend closure extended for two week along loop near i ap news summary at pm edt were now at the militarycheckpoints stage of trump dc takeover trump and putin meet in alaska for highstakes talk sc governor candidate shout nword in cop car on new body cam arrest video trump is trying to make it easier to evict people from public housing gavin newsoms epic trolling is sending republican into fit chicago mayor brandon johnson to trump were not going to cower or bend why is melania trump threatening to sue hunter biden here her story on meeting donald remember the biggest loser a new docuseries explores the harm the show did in the name of health editorial hiding the homeless wont solve the problem biden warned of dc crime crisis in a dems blast trump crackdown dont stop at a stoplight border patrol agent appear outside newsoms event to intimidate trump critic obama call texas gop redistricting an assault on democracy a treaty to end plastic pollution is still out of reach thats not necessarily a bad thing verdict and sentencing for brazil bolsonaro is set for early september protester rally against federal dc police takeover president trump it not too late hear austin tices mother plea to trump year after his disappearance trump armeniaazerbaijan peace effort benefit u company weakens russia newsom cant take a joke this is why trumpputin meeting to take place in alaska today here how to watch live coverage trumpputin meeting to take place in alaska today here how to watch live coverage trump and putin pivotal alaska summit could reshape the war in ukraine rfk jr disavows presidential run pledge loyalty to trump newsom california dems move forward with redistricting effort to counter trumpbacked push in texas more than an overlook reimagining reading pagoda opinion trump and putin to meet at highstakes alaska summit over ending ukraine war poland f fleet is getting a venomous upgrade amazon is selling a pack of polo shirt for just and shopper love their nice athletic cut cold wallet m presale roi projection challenge link whale rally pi network community push best pool and spa construction business loan for bad credit rok financial expands no credit check financing option ahead of season gridstrat unveils gridstrat clarity a business and industryaware ai engine for improved ai search visibility and live market intelligence hoa should release list of board candidate att customer could receive up to in data breach settlement mount charleston property includes ev charger snowmelting roof midyear report show new home sale down unitedhealth group share climb a buffetts berkshire hathaway discloses stake in the insurer nvidia amd see special treatment from trump a chip supply chain face new disruption how alternative financing model are reshaping market best crypto to invest in with explosive roi arctic pablo top the list la vega offer lead paint remediation waster drip repair program summerlin childrens forum award college scholarship how to unplug during vacation and not feel guilty about it denver man called united airline real customer service then got transferred to someone who took k how yahoocom trending coin under could outpace ethereum eth with gain in the coming week gemini ai predicts blockdags daily buyer battle and m miner eclipse bonks drop and ethereums signal could this be the top crypto to buy now even volkswagen is doing horsepower subscription now the drive government shuts down air canada strike that grounded hundred of flight the new york time shareholder reminder faruqi faruqi llp investigates claim on behalf of investor of hims hers health ap business summarybrief at pm edt ap business summarybrief at pm edt government force air canada and flight attendant back to work and into arbitration bessent say he not pushing fed cut just touting model hyundai creta v kia seltos which is better judge block trump ftc retaliation against liberal medium watchdog cnn billionaire investor buy google and this gold miner stock canadian government order arbitration to resolve major air canada flight attendant strike amazon is selling a portable laptop monitor for thats perfect for travel or remote work shareholder reminder faruqi faruqi llp investigates claim on behalf of investor of centene best cannabis business financing expands in a rok financial launch no credit check loan option for licensed operator ap business summarybrief at pm edt government force air canada and flight attendant back to work and into arbitration best smart tv under r with inch above display cashback reward roi how cold wallet is beating xmr hype in s top crypto ranking hottest presale token of august offering shiba inustyle return and x potential shareholder reminder faruqi faruqi llp investigates claim on behalf of investor of petco not enough corvette concept fall hp short of chinese hypercar hype video electrek expandcommunication launch defining intelligent pr billionaire fund manager double down on nvidia and key ai stack play echo investment unit ink million polish real estate sale rillet raise m from andreessen horowitz to replace dumb database with ai accounting that close book in hour amazon is selling a sofa cover for with a boho design that add a lovely touch your money managing cashflow is like planning a summer vacation best contractor business loan for bad credit rok financial responds to growing credit access challenge perplexity ai make surprise pick for best crypto to buy under in here why it neither dogecoin nor sei renewable energy face a looming workforce crisis fly by night catch a luna moth near you trx hype and blockdags ambassador program a strong case for the best crypto platform why strategic change outperforms drastic overhaul shareholder reminder faruqi faruqi llp investigates claim on behalf of investor of sarepta therapeutic larimer county restaurant inspection made aug government force air canada and flight attendant union into arbitration airbnb inc nasdaqabnb share bought by sit investment associate inc bragar eagel squire pc reminds investor of novo altimmune selectquote and spectrum that lawsuit have been filed and encourages investor to contact the firm why your vinyl record is missing song subway and underground garage are untapped geothermal energy source homeowner hit with k invoice past due notice after month of no contact a peco abruptly piled on bill rosen a leading law firm encourages lineage inc investor to secure counsel before jim cramer say dont quit market when it frothy is widespread irrationality a reason to sell down in perfectly rational stock absolutely not jim cramer say dont quit market when it frothy is widespread irrationality a reason to sell down in perfectly rational stock absolutely not cnc investor have opportunity to lead centene corporation security fraud lawsuit canada seek to end airline strike with labor board referral rosen a leading law firm encourages lineage inc investor to secure counsel before important deadline in security class action line how spacepay is changing the conversation about crypto payment usability canadian government intervenes in air canada flight attendant strike a minister order arbitration canadian government intervenes in air canada flight attendant strike a minister order arbitration rosen recognized investor counsel encourages replimune group inc investor to secure rosen recognized investor counsel encourages replimune group inc investor to secure counsel before important deadline in security class action first filed by the firm repl bragar eagel squire pc reminds investor of flywire lockheed martin lineage and sable offshore that lawsuit have been filed and encourages investor to contact the firm amazon take on walmart with new perishable perk it free for million of prime member under a lowered threshold history show that palantir stock monster run is speeding toward an epic crash and it all might come down to detail that no one is talking about history show that palantir stock monster run is speeding toward an epic crash and it all might come down to detail that no one is talking about best crime show on disney to binge this weekend channel rewind grocery price video vault ethereum outpaced cardano but analyst compare the gain from this presale ace and the result are shocking boulder county business hall of fame bert steele made his mark at niwot market larry elder gerrymandering obama v obama tesla raise model sx price but bundle full selfdriving and brings back free supercharging airbnb inc nasdaqabnb stock holding boosted by vestmark advisory solution inc best no credit check business financing demand rise a rok financial expands access to small business loan in greenlight capital announces q trade new position in flr and vsco exit from vtrs a startup is tapping underground parking garage for clean energy employer need help managing worker who are taking second job waste management inc nysewm share acquired by horizon investment llc highgain token to watch a early buyer target the best crypto coin with x potential in moonshot magax redefining viral content monetization with loomint ai walmart is selling a roomy buffet storage cabinet for only letter another idea to save city money ford motor co ha a much bigger issue than tariff at it door target and ulta mutually agree not to renew partnership launched in boston family office llc boost stock position in alphabet inc nasdaqgoog forvis mazars wealth advisor llc ha million position in walmart inc nysewmt top crypto presales for dominating the market with huge roi potential procter gamble company the nysepg share acquired by schneider down wealth management advisor lp business briefcase explore schuylkill tout county attraction developer seek to build condo complex for downsizing senior social security trustee say the program can pay all benefit until it might be too optimistic the redneck recession trump offer putin zelensky contrasting approach south park actress trump supporter janeshia adamsginyard happy show is taking him on west virginia sends hundred of national guard member to washington at trump team request trump rate putin meeting a out of because they got along great trump rate putin meeting a out of because they got along great dear abby woman wish to acknowledge kindness after tragic loss about those judge who take side reader commentary key takeaway from the trumpputin summit three republicanled state to send hundred of national guard troop to washington west virginia sends hundred of national guard member to washington at trump team request willis seeking fifth term a wetumpka mayor trump told zelenskyy after summit that putin want more of ukraine reuters source say people in ice custody complain of long wait for deportation state delay new contract with extra accountability measure for pa elder abuse agency lock them all up kash patel declassifies evidence of doj coverup of clinton foundation corruption robert f kennedy jr declares loyalty to trump rule out a presidential bid shareholder reminder faruqi faruqi llp investigates claim on behalf of investor of hims hers health government force air canada and flight attendant back to work and into arbitration court could decide when search pas the smell test the trumpputin summit wa a win for russia keep racial lunacy out of air traffic control before more plane crash west virginia governor is sending to national guard member to washington at the trump administration request judge block trump ftc retaliation against liberal medium watchdog photo of the putintrump summit in alaska immigration raid and fear of deportation push patient into telehealth the fight is on how redistricting could unfold in entangled state people are completely losing it over this maga voter rant telling donald trump to stop living in the past and demanding to see some result hawk fire now contained in acton trump warned by top senate democrat to rethink advanced ai chip sale to china summit put putin back on the global stage and trump echo a kremlin position summit put putin back on the global stage and trump echo a kremlin position summit put putin back on the global stage and trump echo a kremlin position mckees rock teen charged with attempted homicide in shooting at east liberty garland park mckees rock teen charged with attempted homicide in east liberty park shooting clip of the week texas democrat outfoxes and rfk jr get roasted the virginiaclass submarine is stepping up a the ssnx face delay opinion california is no longer a model for transparent government these two plane could fill the gap left by the a warthog shareholder reminder faruqi faruqi llp investigates claim on behalf of investor of sarepta therapeutic credit when it due bill maher give trump high mark for tariff nato and being antiwar watch democrat will soon fall to third place beto orourke say his political org donated more than m during first special session focus on mental health strike at root of gun violence reader commentary developman editorial cartoon ten way to force putin back to the bargaining table bill maher rip hunter biden after melania trump threatens to sue disgraced exfirst son campbell pb nc defunded not defeated west virginia national guard to help trump restore safety in nation capital folk in dc already feeling safer praise trump for crime crackdown rosen recognized investor counsel encourages replimune group inc investor to secure counsel before important deadline in security class action first filed by the firm repl bragar eagel squire pc reminds investor of flywire lockheed martin lineage and sable offshore that lawsuit have been filed and encourages investor to contact the firm uscis to more heavily scrutinize good moral character citizenship requirement trump administration to more heavily scrutinize good moral character requirement for u citizenship hurricane erin path spaghetti model tracker will category storm impact florida vestmark advisory solution inc purchase share of airbnb inc nasdaqabnb leader should be candid on crime not posture editorial horizon investment llc ha million stake in waste management inc nysewm letter we need access to boca beach putin emerges from the alaska summit with increased stature and trump echoing a kremlin position alphabet inc nasdaqgoog is progressive investment management corp th largest position letter trump plan for epstein file is confusing letter minority should not govern our nation notice a theme to trump planned takeover of city these black mayor do cheering for war jessica tarlov gloat that trump doesnt need to worry about the nobel prize in world war ii the soviet union built a flying tank letter we need to make our voice heard your view keep outside money and politics out of our school do cuomo and adam secretly want mamdani to win could the hellhound drone be a contender for the golden dome put pressure on putin mr president it the only way to end this war alliance wealth advisor llc ut acquires share of salesforce inc nysecrm the home depot inc nysehd share purchased by alliance wealth advisor llc ut zoetis inc nysezts share acquired by vestmark advisory solution inc armstrong williams the silent majority are applauding staff commentary adobe inc nasdaqadbe share bought by tocqueville asset management lp cnn anchor jake tapper caught on hot mic during trumpputin coverage give me my show ap news summary at am edt letter let all be peakbusters letter there such a thing a integrity menora mivtachim holding ltd boost holding in invesco qqq nasdaqqqq cookson peirce co inc buy new position in invesco qqq nasdaqqqq cookson peirce co inc ha holding in waste management inc nysewm american century company inc raise position in targa resource inc nysetrgp letter the dnc should take lesson from nd demnpl party callan family office llc decrease stake in meta platform inc nasdaqmeta rossmore private capital boost position in meta platform inc nasdaqmeta liontrust investment partner llp sell share of water corporation nysewat u gut mention of indigenous right climate change from human right report reddit urgent deadline alert bragar eagel squire pc reminds reddit investor of the august th deadline for the class action suit alphabet inc nasdaqgoog is boston family office llcs th largest position eli lilly and company nyselly position decreased by roanoke asset management corp ny fox run management llc increase position in thermo fisher scientific inc nysetmo vanguard sp etf nysearcavoo share purchased by twin focus capital partner llc los gatos resident upset over reckless ebike rider letter our right are being degraded schneider down wealth management advisor lp ha holding in procter gamble company the nysepg american express company nyseaxp share sold by lord abbett co llc gulf international bank uk ltd reduces stock holding in equity residential nyseeqr mcdonalds corporation nysemcd share bought by vestmark advisory solution inc aberdeen group plc grows stock position in kinder morgan inc nysekmi american drown in a hour span at same bahamas resort police abc news magacoin finance price prediction analyst see x gain before xrp etf approval mother sue florida dairy farm claiming she lost unborn baby after toddler got sick from drinking raw milk cbs news cold wallet raise m while bitcoin cash halving spark momentum xrp whale move push speculation higher ag gov asked him to probe chinese company land buy in nashua how high dogecoin could go if it lead bull run and why new doge rival under could skyrocket this hightech bra is flying off shelf thanks to taylor swift ap business summarybrief at pm edt air canada suspends restart plan after flight attendant union defies return to work order bnb price reach new alltime high nexchain presale earns best crypto to buy status for vega luxury home market buck housing slump uber freight ceo lior ron leaf b unit for waabis aifirst driverless truck most positioned to lead the transformation state rethink a longheld practice of setting speed limit based on how fast driver travel fatal explosion at u steel plant raise question about it future despite heavy investment fatal explosion at u steel plant raise question about it future despite heavy investment ap business summarybrief at pm edt air canada suspends restart plan after flight attendant union defies return to work order air canada suspends restart plan after flight attendant union defies return to work order chesapeake native work with hampton innovation center to grow farm system gayle cook cook group cofounder and preservationist dy at air canada suspends restart plan after union defies return to work order peter lynch stock market ha been the best place to be but if you need money in or year you shouldnt be buying stock noise v boat best budget neckband in rising food price are likely to continue thanks to trump trade war rolling stone seven mile advisory increase stake in well fargo company nysewfc air canada suspends restart plan after union defies return to work order duolingos ai shift spark backlash amid revenue growth klc investor notice robbins geller rudman dowd llp file class action lawsuit against kindercare learning company inc and announces opportunity for investor with substantial loss to lead class action lawsuit foyston gordon payne inc ha million stock holding in unitedhealth group incorporated nyseunh air canada suspends plan to restart operation on sunday after flight attendant union defies return to work order air canada suspends plan to restart operation on sunday after flight attendant union defies return to work order here the share of gold or crypto ray dalio say investor should hold if youd invested in shib year ago here how much youd have today freestone capital holding llc sell share of chevron corporation nysecvx bolivia presidential vote headed for firstever runoff between centrist rightwing candidate reason you should switch to an electric vehicle in retirement gobankingrates cold wallet lead with proven utility before token launch cardano target ethena struggle at cardanos ada modest jump fails to retain tired holder investor shift to new coin under fo'
$ws.Cells.Item($newRow, 10).Value = 'strategy mstr underperforming bitcoin whats driving divergence online casino company acquires spotlightvegas medium advisory federal government housing announcement stratford cold wallet among highest trending crypto raised xrp stay bullish dot price eye colorado springsbased auria million contract watch trump putin share warm greeting handshake alaska summit kick rddt deadline monday rosen ranked law firm encourages reddit investor air canada shut completely unless airline flight attendant reach deal air canada shut completely unless airline flight attendant reach deal advocate mark passage ct towing reform law jetblue add florida flight msnbc mouthpiece claim climate change dangerous crime dc caitong international asset management co acquires share visa nysev rddt deadline monday rosen ranked law firm encourages reddit investor loss excess k secure counsel important august deadline security class action rddt rosen recognized investor counsel encourages biohaven investor secure counsel target ulta mutually agree renew partnership launched rosen recognized investor counsel encourages biohaven investor secure counsel important deadline security class action bhvn resident biggest tax benefit trump big bill trump open business russia tricky tariff sale slump hit deere analyst greener pasture ahead duke energy seek extend operation another year bad creek supporting unprecedented growth carolina inventhelp inventor develops new laserbased tool plumber tro government canada provides early decision york factory nation shilling aerodrome project manitoba saks hyperpersonalized homepage boost revenue per visitor monthly payment fall month low report trump new usrussia trade deal war settledahead putin meeting update job coming dearborn county federal grand jury indicts new orleans mayor corruption probe little pepe lilpepe price forecast unrealistic longterm target gamblingcom analyst slash forecast following q forward industry recovers q revenue decline earnings miss halle berry busy cooking cleaning mothering picture perfect response ex upsidedown car loan breaking sandisk analyst share view analyst rating restaurant brand intl assessing netstreit insight financial analyst unpacking latest option trading trend centrus energy summer travel isnt easy used airline analyst revise forecast deere q spotlight conocophillips analyzing surge option activity demystifying amcor insight analyst review analyst saying nnn reit stock farleys model moment bold ev strategy risky policy era merchant corporation announces cash dividend pe ratio insight targa resource analyst saying broadstone net lease stock inside backyard housing boom california north river marketplace celebrates year selling food resident jollibee global girl katseye ignite partnership create joyful new connection north american fan presale blockdag academy put stellar bittensor rearview best crypto buy sunrun run stock soar fargo nearly double price target liverpool sign teenage italian defender giovanni leoni liverpool sign teenage italian defender giovanni leoni cfo must take ownership managing culture risk auditboard analyst increase forecast advance auto part following upbeat q earnings novo nordisks lead fade eli lilly course outsized gain analyst nearly year legendary japanese sport car unexpected twist gear patrol fact focus taxpayer receive new stimulus summer fact focus taxpayer receive new stimulus summer fact focus taxpayer receive new stimulus summer fact focus taxpayer receive new stimulus summer fact focus taxpayer receive new stimulus summer got microplastic pollution anxiety chris hemsworth cure analyst revise forecast applied industrial technology following q earnings layer brett pepe forecast fight ethereums leading meme spot begin californian afford v nationally strike loom air canada flight attendant reject call arbitration new york air canada flight attendant reject call arbitration unitedhealths best day since intels best weekly gain decade whats moving market friday special livery package gt celebrates anniversary andial trump america ceo dogecoin news today dogecoin whale transfer binance trigger volatility concern senator call probe meta news report ai conversation child sharplink share drop q revenue decline large noncash eth impairment retail sale rose july shopper stepped purchase ahead tariff trump intel great tsmc fed end specialized crypto bank oversight folding normal supervision analyst boost forecast dillards following strong q earnings new toll brother model open station sunnyvale california burger bending budget blame climate change best memecoin buy pepeto next bnb memecoins turn investment million agco order immediate suspension north bay max multi mart lottery alcohol licence response human trafficking charge washington dc trump administration court hearing police takeover lawsuit stanley firstever warehouse sale deal employee longtime newport news travel agency year experience combined industrials stock whale activity today session popular bakery hoping avoid closure zoning complaint samsung giving away tb ssd odyssey g monitor deepwater asset management llc million stock holding broadcom nasdaqavgo billionaire family office join berkshire unitedhealth bet yahoo finance marlton partner file preliminary proxy statement related election director degree capital board director lego newest tuxedo cat purrrfect addition cat lover jeff bezos lauren sanchez suffer rocky start marriage devastating update yahoocom jeff bezos lauren sanchez suffer rocky start marriage devastating update sheknows ai human expertise changing wealth management service target ulta beauty ending store partnership target ulta beauty ending store partnership solana whale activity resume stellars rally cold wallet viral raised sol price swing push solana investor layer brett staking reward coding apps pick best perrone unlike columbia trump attack ucla aimed taxpayer money uhoh trump rizz among youth voted republican pitch statewide office senate illinois dems blueprint national fight trump texas gop chairman call illinois winnable purple house minority keynote speaker annual democratic county chair brunch editorial hiding homeless wont solve problem shribman future pat buchanan letter trump administration mistreats federal worker federal grand jury indicts new orleans mayor latoya cantrell corruption probe new orleans mayor indicted allegation hide relationship bodyguard federal grand jury indicts new orleans mayor latoya cantrell corruption probe federal grand jury indicts new orleans mayor latoya cantrell corruption probe federal grand jury indicts new orleans mayor latoya cantrell corruption probe federal grand jury indicts new orleans mayor latoya cantrell corruption probe federal grand jury indicts new orleans mayor latoya cantrell corruption probe new orleans mayor latoya cantrell indicted federal grand jury following lengthy corruption investigation texas gov abbott call immediate special session redistricting bad diet med exercise look inside latest america healthy report trump arrives alaska military base summit putin russiaukraine war trump arrives alaska military base summit putin russiaukraine war trump arrives alaska military base summit putin russiaukraine war bad diet med exercise look inside latest america healthy report trump open business russia tricky antivaxxers making excuse cdc shooter putin tease new nuke deal ahead alaska meeting trump talking past useu dispute free speech set escalate trump new usrussia trade deal war settledahead putin meeting update trump arrives anchorage alaska federal grand jury indicts new orleans mayor corruption probe adyens tariff trouble portend payment company marjorie taylor greene buying unitedhealth stock along fellow investor buffett burry new orleans mayor indicted allegation hide relationship bodyguard halle berry busy cooking cleaning mothering picture perfect response ex audience erupts rep jasmine crockett town hall amid verbal attack audience erupts rep jasmine crockett town hall amid verbal attack america national park scenic theyre sacred created cost native american cornyn paxton dead heat senate primary latest poll plan open new detention center rural colorado report plan open new detention center rural colorado report automatic gratuity eligible new trump tax break tip ai changing music business switzerland scale f order following tariff threat republican view supporting ukraine swing dramatically ahead trump meeting putin new poll kiss criticized trump prior kennedy center honor fact focus taxpayer receive new stimulus summer fact focus taxpayer receive new stimulus summer fact focus taxpayer receive new stimulus summer fact focus taxpayer receive new stimulus summer fact focus taxpayer receive new stimulus summer fact focus taxpayer receive new stimulus summer fact focus taxpayer receive new stimulus summer fact focus taxpayer receive new stimulus summer white house fire trump dc police crime crackdown placed crosshairs dems witnessing deportation devastating child dhs noem temporarily staying military housing facing vicious doxxing death threat texas house speaker expects enough democrat quorum texas house speaker expects enough democrat quorum texas house speaker expects enough democrat quorum texas house speaker expects enough democrat quorum texas house speaker expects enough democrat quorum texas house democrat accept skokie invite north suburban mosque news summary pm la nina weather california trump launched new assault dc rule brazenly illegal consumer urge fhfa mix crypto mortgage washington dc trump administration court hearing police takeover lawsuit newsom end trump presidency reveals plan bypass ca redistricting commission trump putin pivotal alaska summit reshape war ukraine trump mull lower refugee cap space white south african hawley open probe meta report ai romantic exchange minor septum urge rider start planning upcoming service cut begin aug weekend closure extended along loop news summary pm militarycheckpoints stage trump dc takeover trump putin meet alaska highstakes sc governor candidate shout nword cop car new body cam arrest video trump easier evict people public housing gavin newsoms epic trolling sending republican fit chicago mayor brandon johnson trump cower bend melania trump threatening sue hunter biden story meeting donald remember biggest loser new docuseries explores harm name health editorial hiding homeless wont solve problem biden warned dc crime crisis dems blast trump crackdown dont stop stoplight border patrol agent appear outside newsoms event intimidate trump critic obama call texas gop redistricting assault democracy treaty end plastic pollution reach thats necessarily bad verdict sentencing brazil bolsonaro set early september protester rally federal dc police takeover president trump late hear austin tices mother plea trump year disappearance trump armeniaazerbaijan peace effort benefit company weakens russia newsom cant take joke trumpputin meeting take place alaska today watch coverage trumpputin meeting take place alaska today watch coverage trump putin pivotal alaska summit reshape war ukraine rfk jr disavows presidential run pledge loyalty trump newsom california dems move forward redistricting effort counter trumpbacked push texas overlook reimagining reading pagoda opinion trump putin meet highstakes alaska summit ending ukraine war poland f fleet getting venomous upgrade amazon selling pack polo shirt shopper love nice athletic cut cold wallet presale roi projection challenge link whale rally pi network community push best pool spa construction business loan bad credit rok financial expands credit financing option ahead season gridstrat unveils gridstrat clarity business industryaware ai engine improved ai search visibility market intelligence hoa release list board candidate att customer receive data breach settlement mount charleston property includes ev charger snowmelting roof midyear report new sale unitedhealth share climb buffetts berkshire hathaway discloses stake insurer nvidia amd special treatment trump chip supply chain face new disruption alternative financing model reshaping market best crypto invest explosive roi arctic pablo list la vega offer lead paint remediation waster drip repair program summerlin childrens forum award college scholarship unplug vacation feel guilty denver man united airline real customer service got transferred someone took k yahoocom trending coin outpace ethereum eth gain coming gemini ai predicts blockdags buyer battle miner eclipse bonks drop ethereums signal crypto buy volkswagen horsepower subscription drive government shuts air canada strike grounded flight new york shareholder reminder faruqi faruqi llp investigates claim behalf investor hims health business summarybrief pm business summarybrief pm government force air canada flight attendant work arbitration bessent pushing fed cut touting model hyundai creta v kia seltos better judge block trump ftc retaliation liberal medium watchdog cnn billionaire investor buy google gold miner stock canadian government order arbitration resolve major air canada flight attendant strike amazon selling portable laptop monitor thats perfect travel remote work shareholder reminder faruqi faruqi llp investigates claim behalf investor centene best cannabis business financing expands rok financial launch credit loan option licensed operator business summarybrief pm government force air canada flight attendant work arbitration best smart tv r inch display cashback reward roi cold wallet beating xmr hype crypto ranking hottest presale token august offering shiba inustyle x potential shareholder reminder faruqi faruqi llp investigates claim behalf investor petco enough corvette concept fall hp short chinese hypercar hype video electrek expandcommunication launch defining intelligent pr billionaire fund manager double nvidia ai stack play echo investment unit ink million polish real estate sale rillet raise andreessen horowitz replace dumb database ai accounting close book hour amazon selling sofa cover boho design add lovely touch money managing cashflow planning summer vacation best contractor business loan bad credit rok financial responds growing credit access challenge perplexity ai surprise pick best crypto buy neither dogecoin sei renewable energy face looming workforce crisis fly night catch luna moth trx hype blockdags ambassador program strong case best crypto platform strategic change outperforms drastic overhaul shareholder reminder faruqi faruqi llp investigates claim behalf investor sarepta therapeutic larimer county restaurant inspection aug government force air canada flight attendant union arbitration airbnb nasdaqabnb share bought sit investment associate bragar eagel squire pc reminds investor novo altimmune selectquote spectrum lawsuit filed encourages investor contact firm vinyl record missing song subway underground garage untapped geothermal energy source homeowner hit k invoice past due notice month contact peco abruptly piled bill rosen leading law firm encourages lineage investor secure counsel jim cramer dont quit market frothy widespread irrationality sell perfectly rational stock absolutely jim cramer dont quit market frothy widespread irrationality sell perfectly rational stock absolutely cnc investor opportunity lead centene corporation security fraud lawsuit canada seek end airline strike labor board referral rosen leading law firm encourages lineage investor secure counsel important deadline security class action line spacepay changing conversation crypto payment usability canadian government intervenes air canada flight attendant strike minister order arbitration canadian government intervenes air canada flight attendant strike minister order arbitration rosen recognized investor counsel encourages replimune investor secure rosen recognized investor counsel encourages replimune investor secure counsel important deadline security class action filed firm repl bragar eagel squire pc reminds investor flywire lockheed martin lineage sable offshore lawsuit filed encourages investor contact firm amazon take walmart new perishable perk free million prime lowered threshold history palantir stock monster run speeding epic crash come detail talking history palantir stock monster run speeding epic crash come detail talking best crime disney binge weekend channel rewind grocery price video vault ethereum outpaced cardano analyst compare gain presale ace shocking boulder county business hall fame bert steele mark niwot market larry elder gerrymandering obama v obama tesla raise model sx price bundle full selfdriving brings free supercharging airbnb nasdaqabnb stock holding boosted vestmark advisory solution best credit business financing demand rise rok financial expands access small business loan greenlight capital announces q trade new flr vsco exit vtrs startup tapping underground parking garage clean energy employer help managing worker taking job waste management nysewm share acquired horizon investment llc highgain token watch early buyer target best crypto coin x potential moonshot magax redefining viral content monetization loomint ai walmart selling roomy buffet storage cabinet letter another idea save money ford co bigger issue tariff door target ulta mutually agree renew partnership launched boston family office llc boost stock alphabet nasdaqgoog forvis mazars wealth advisor llc million walmart nysewmt crypto presales dominating market roi potential procter gamble company nysepg share acquired schneider wealth management advisor lp business briefcase explore schuylkill tout county attraction developer seek build condo complex downsizing senior social security trustee program pay benefit optimistic redneck recession trump offer putin zelensky contrasting approach south park actress trump supporter janeshia adamsginyard happy taking west virginia sends national guard washington trump team request trump rate putin meeting got along great trump rate putin meeting got along great dear abby woman wish acknowledge kindness tragic loss judge take reader commentary takeaway trumpputin summit republicanled send national guard troop washington west virginia sends national guard washington trump team request willis seeking fifth term wetumpka mayor trump zelenskyy summit putin ukraine reuters source people custody complain wait deportation delay new contract extra accountability measure pa elder abuse agency lock kash patel declassifies evidence doj coverup clinton foundation corruption robert f kennedy jr declares loyalty trump rule presidential bid shareholder reminder faruqi faruqi llp investigates claim behalf investor hims health government force air canada flight attendant work arbitration court decide search pa smell test trumpputin summit russia racial lunacy air traffic control plane crash west virginia governor sending national guard washington trump administration request judge block trump ftc retaliation liberal medium watchdog photo putintrump summit alaska immigration raid fear deportation push patient telehealth fight redistricting unfold entangled people completely losing maga voter rant telling donald trump stop living past demanding hawk fire contained acton trump warned senate democrat rethink advanced ai chip sale china summit put putin global stage trump echo kremlin summit put putin global stage trump echo kremlin summit put putin global stage trump echo kremlin mckees rock teen charged attempted homicide shooting east liberty garland park mckees rock teen charged attempted homicide east liberty park shooting clip texas democrat outfoxes rfk jr roasted virginiaclass submarine stepping ssnx face delay opinion california longer model transparent government plane fill gap left warthog shareholder reminder faruqi faruqi llp investigates claim behalf investor sarepta therapeutic credit due bill maher give trump high mark tariff nato antiwar watch democrat soon fall third place beto orourke political org donated special session focus mental health strike root gun violence reader commentary developman editorial cartoon force putin bargaining table bill maher rip hunter biden melania trump threatens sue disgraced exfirst son campbell pb nc defunded defeated west virginia national guard help trump restore safety nation capital folk dc feeling safer praise trump crime crackdown rosen recognized investor counsel encourages replimune investor secure counsel important deadline security class action filed firm repl bragar eagel squire pc reminds investor flywire lockheed martin lineage sable offshore lawsuit filed encourages investor contact firm uscis heavily scrutinize good moral character citizenship requirement trump administration heavily scrutinize good moral character requirement citizenship hurricane erin path spaghetti model tracker category storm impact florida vestmark advisory solution purchase share airbnb nasdaqabnb candid crime posture editorial horizon investment llc million stake waste management nysewm letter access boca beach putin emerges alaska summit increased stature trump echoing kremlin alphabet nasdaqgoog progressive investment management corp largest letter trump plan epstein file confusing letter minority govern nation notice theme trump planned takeover black mayor cheering war jessica tarlov gloat trump doesnt worry nobel prize world war ii soviet union built flying tank letter voice heard view outside money politics school cuomo adam secretly mamdani hellhound drone contender golden dome put pressure putin mr president end war alliance wealth advisor llc ut acquires share salesforce nysecrm depot nysehd share purchased alliance wealth advisor llc ut zoetis nysezts share acquired vestmark advisory solution armstrong williams silent majority applauding staff commentary adobe nasdaqadbe share bought tocqueville asset management lp cnn anchor jake tapper caught hot mic trumpputin coverage give news summary letter peakbusters letter integrity menora mivtachim holding boost holding invesco qqq nasdaqqqq cookson peirce co buy new invesco qqq nasdaqqqq cookson peirce co holding waste management nysewm american century company raise targa resource nysetrgp letter dnc take lesson nd demnpl party callan family office llc decrease stake meta platform nasdaqmeta rossmore private capital boost meta platform nasdaqmeta liontrust investment partner llp sell share water corporation nysewat gut mention indigenous climate change human report reddit urgent deadline alert bragar eagel squire pc reminds reddit investor august deadline class action suit alphabet nasdaqgoog boston family office llcs largest eli lilly company nyselly decreased roanoke asset management corp ny fox run management llc increase thermo fisher scientific nysetmo vanguard sp etf nysearcavoo share purchased twin focus capital partner llc los gatos resident upset reckless ebike rider letter degraded schneider wealth management advisor lp holding procter gamble company nysepg american express company nyseaxp share sold lord abbett co llc gulf international bank uk reduces stock holding equity residential nyseeqr mcdonalds corporation nysemcd share bought vestmark advisory solution aberdeen plc grows stock kinder morgan nysekmi american drown hour span bahamas resort police abc news magacoin finance price prediction analyst x gain xrp etf approval mother sue florida dairy claiming unborn baby toddler got sick drinking raw milk cbs news cold wallet raise bitcoin cash halving spark momentum xrp whale move push speculation higher ag gov probe chinese company land buy nashua high dogecoin lead bull run new doge rival skyrocket hightech bra flying shelf taylor swift business summarybrief pm air canada suspends restart plan flight attendant union defies work order bnb price reach new alltime high nexchain presale earns best crypto buy status vega luxury market buck housing slump uber freight ceo lior ron leaf b unit waabis aifirst driverless truck positioned lead transformation rethink longheld practice setting speed limit fast driver travel fatal explosion steel plant raise future heavy investment fatal explosion steel plant raise future heavy investment business summarybrief pm air canada suspends restart plan flight attendant union defies work order air canada suspends restart plan flight attendant union defies work order chesapeake native work hampton innovation center grow system gayle cook cook cofounder preservationist dy air canada suspends restart plan union defies work order peter lynch stock market best place money year shouldnt buying stock noise v boat best budget neckband rising food price likely continue trump trade war rolling stone mile advisory increase stake fargo company nysewfc air canada suspends restart plan union defies work order duolingos ai shift spark backlash amid revenue growth klc investor notice robbins geller rudman dowd llp file class action lawsuit kindercare learning company announces opportunity investor substantial loss lead class action lawsuit foyston gordon payne million stock holding unitedhealth incorporated nyseunh air canada suspends plan restart operation sunday flight attendant union defies work order air canada suspends plan restart operation sunday flight attendant union defies work order share gold crypto ray dalio investor youd invested shib year ago youd today freestone capital holding llc sell share chevron corporation nysecvx bolivia presidential vote headed firstever runoff centrist rightwing candidate switch electric vehicle retirement gobankingrates cold wallet lead proven utility token launch cardano target ethena struggle cardanos ada modest jump fails retain tired holder investor shift new coin fo'
